$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (52) down into the four
# new rows so the new date cells (column A) pick up the same date number
# format (style) as the rest of the table.
$ws.Range("A52:G52").Copy() | Out-Null
$ws.Range("A53:G56").PasteSpecial(-4122) | Out-Null

# New daily log entries (dates 3/12/2021 - 3/15/2021, serials 44267-44270)
$data = @(
    @(44267, 0, 0, 0, 0, 0, 0),
    @(44268, 0, 0, 0, 0, 0, 0),
    @(44269, 0, 0, 0, 0, 0, 0),
    @(44270, 60, 45, 0, 30, 100, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 53 + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}

# Scroll/select to match the author's final view position.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("J52").Select() | Out-Null
